$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Rename the header row cells: "<Name>_old" -> "<Name>_FV2404",
#    "<Name>_new" -> "<Name>_FV2410". Column K ("diff") is left untouched.
$oldHeaders = @(
    "Segmentname_old", "Segmentgruppe_old", "Segment_old", "Datenelement_old",
    "Segment ID_old", "Code_old", "Qualifier_old", "Beschreibung_old",
    "Bedingungsausdruck_old", "Bedingung_old"
)
$newHeaders = @(
    "Segmentname_new", "Segmentgruppe_new", "Segment_new", "Datenelement_new",
    "Segment ID_new", "Code_new", "Qualifier_new", "Beschreibung_new",
    "Bedingungsausdruck_new", "Bedingung_new"
)

for ($i = 0; $i -lt $oldHeaders.Count; $i++) {
    $col = $i + 1
    $label = $oldHeaders[$i] -replace "_old$", "_FV2404"
    $ws.Cells.Item(1, $col).Value = $label
}

for ($i = 0; $i -lt $newHeaders.Count; $i++) {
    $col = $i + 12
    $label = $newHeaders[$i] -replace "_new$", "_FV2410"
    $ws.Cells.Item(1, $col).Value = $label
}

# 2) Turn the used range into an Excel Table ("Table1") so the header row
#    carries filter buttons and the workbook gets xl/tables/table1.xml.
#    Stash the existing header formatting first and restore it afterwards
#    (ListObjects.Add otherwise bakes the header's current format into a
#    brand-new dxf / headerRowDxfId instead of reusing the existing style).
$headerRange = $ws.Range("A1:U1")
$scratch = $ws.Range("A200:U200")
$headerRange.Copy()
$scratch.PasteSpecial(-4122)  # xlPasteFormats
$headerRange.ClearFormats()

$rng = $ws.Range("A1:U84")
$tbl = $ws.ListObjects.Add(1, $rng, [System.Reflection.Missing]::Value, 1)
$tbl.Name = "Table1"

$scratch.Copy()
$headerRange.PasteSpecial(-4122)  # xlPasteFormats
$scratch.ClearFormats()
$scratch.ClearContents()

# 3) Freeze the header row (split below row 1, top-left of the scrollable
#    area is A2).
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
